$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 189: first session of 2023-09-29 (17:30 - 21:30)
$ws.Range("A189").Value = "2023-09-29"
$ws.Range("B189").Value = "17:30"
$ws.Range("C189").Value = "21:30"
$ws.Range("D189").Value = "4h 00m"
$ws.Range("E189").Value = "#python"
$ws.Range("F189").Value = "nwtimetrackingmanager v1.0.0"
$ws.Range("G189").Value = "'True"
$ws.Range("H189").Value = "'False"
$ws.Range("I189").Formula = "=YEAR(A189)"
$ws.Range("J189").Formula = "=MONTH(A189)"

# Row 190: second session of 2023-09-29 (22:45 - 03:15, overnight)
$ws.Range("A190").Value = "2023-09-29"
$ws.Range("B190").Value = "22:45"
$ws.Range("C190").Value = "03:15"
$ws.Range("D190").Value = "4h 30m"
$ws.Range("E190").Value = "#python"
$ws.Range("F190").Value = "nwtimetrackingmanager v1.0.0"
$ws.Range("G190").Value = "'True"
$ws.Range("H190").Value = "'False"
$ws.Range("I190").Formula = "=YEAR(A190)"
$ws.Range("J190").Formula = "=MONTH(A190)"

# Move selection / active cell like the recorded edit
$ws.Range("E190").Select() | Out-Null
